$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 356; this shifts old rows 356..383 down to 357..384,
# matching the weekly data-refresh described in the commit message.
$ws.Rows.Item(356).Insert()

# Populate the newly inserted row 356 with this week's new record.
$ws.Range("A356").Value = 11
$ws.Range("B356").Value = "Vega Monumental Concepción"
$ws.Range("C356").Value = "Bíobío"
$ws.Range("D356").Value = 44931
$ws.Range("E356").Value = 8
$ws.Range("F356").Value = 100112017
$ws.Range("G356").Value = "Apio"
$ws.Range("H356").Value = "Americana (o)"
$ws.Range("I356").Value = "Primera"
$ws.Range("J356").Value = 350
$ws.Range("K356").Value = 9000
$ws.Range("L356").Value = 9500
$ws.Range("M356").Value = 9286
$ws.Range("N356").Value = "`$/docena de matas"
$ws.Range("O356").Value = "Región de Coquimbo"
$ws.Range("P356").Value = 1548
$ws.Range("Q356").Value = 6
$ws.Range("R356").Value = "Hortaliza"
